$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# This edit adds two new HKL rows ("Holden" and "Rizzie Spiral") right after
# "Spiral5", renames "Thomas Hex" to "Matthies Hex", and appends two more new
# HKL rows ("Michael-CCHex" and "Michael-SNHex") at the bottom, re-using the
# simulation output that used to sit in the last two rows. All of this is
# expressed as direct cell writes driven by the already re-run simulation
# output, so every cell ends up holding its final value without relying on
# row-shifting operations.
# ---------------------------------------------------------------------------

# Extend the A column formatting (bold, centered, thin border) used for the
# index column down through the new rows 4, 5, 30 and 31 by copying the
# existing format from A3.
$ws.Range("A3").Copy() | Out-Null
$ws.Range("A4:A5").PasteSpecial(-4122) | Out-Null
$ws.Range("A30:A31").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Row 4: HKL index 2 = "Holden"
$ws.Cells.Item(4,1).Value = 2
$ws.Cells.Item(4,2).Value = "Holden"
$r4 = New-Object 'object[,]' 1,18
$r4[0,0] = 0.999669715406289
$r4[0,1] = 0.7813224481860843
$r4[0,2] = 0.7813224481860843
$r4[0,3] = 0.9018466284670321
$r4[0,4] = 0.9018466284670321
$r4[0,5] = 0.809922717544755
$r4[0,6] = 3.303499476935772
$r4[0,7] = 0.9487499753379449
$r4[0,8] = 0.9018466284670321
$r4[0,9] = 0.999669715406289
$r4[0,10] = 0.8904960817961867
$r4[0,11] = 0.8904960817961867
$r4[0,12] = 0.8636382937123761
$r4[0,13] = 0.8942795973531351
$r4[0,14] = 0.8942795973531351
$r4[0,15] = 0.8961713551316094
$r4[0,16] = 0.8961713551316094
$r4[0,17] = 1.29083516031298
$ws.Range("C4:T4").Value = $r4

# Row 5: HKL index 3 = "Rizzie Spiral"
$ws.Cells.Item(5,1).Value = 3
$ws.Cells.Item(5,2).Value = "Rizzie Spiral"
$r5 = New-Object 'object[,]' 1,18
$r5[0,0] = 1.946201263933978
$r5[0,1] = 3.769077610236998
$r5[0,2] = 3.769077610236998
$r5[0,3] = 1.882590041464107
$r5[0,4] = 1.882590041464107
$r5[0,5] = 0.03548979216606776
$r5[0,6] = 3.919570077975332
$r5[0,7] = 0.1425334601701432
$r5[0,8] = 1.882590041464107
$r5[0,9] = 1.946201263933978
$r5[0,10] = 2.857639437085489
$r5[0,11] = 2.857639437085489
$r5[0,12] = 1.916922888779015
$r5[0,13] = 2.532622971878361
$r5[0,14] = 2.532622971878361
$r5[0,15] = 2.370114739274797
$r5[0,16] = 2.370114739274797
$r5[0,17] = 1.949243707657771
$ws.Range("C5:T5").Value = $r5

# Row 6: HKL index 4 = "RotRing OmegaMax-90"
$ws.Cells.Item(6,1).Value = 4
$ws.Cells.Item(6,2).Value = "RotRing OmegaMax-90"
$r6 = New-Object 'object[,]' 1,18
$r6[0,0] = 0.8840667873860256
$r6[0,1] = 0.8086297637918995
$r6[0,2] = 0.8086297637918995
$r6[0,3] = 1.891452413364914
$r6[0,4] = 1.891452413364914
$r6[0,5] = 0.8630254139768738
$r6[0,6] = 0.933098027846971
$r6[0,7] = 0.9716474003225971
$r6[0,8] = 1.891452413364914
$r6[0,9] = 0.8840667873860256
$r6[0,10] = 0.8463482755889626
$r6[0,11] = 0.8463482755889626
$r6[0,12] = 0.8519073217182663
$r6[0,13] = 1.19471632151428
$r6[0,14] = 1.19471632151428
$r6[0,15] = 1.368900344476938
$r6[0,16] = 1.368900344476938
$r6[0,17] = 1.05865330111488
$ws.Range("C6:T6").Value = $r6

# Row 7: HKL index 5 = "Equal Angle"
$ws.Cells.Item(7,1).Value = 5
$ws.Cells.Item(7,2).Value = "Equal Angle"
$r7 = New-Object 'object[,]' 1,18
$r7[0,0] = 1.044529099493631
$r7[0,1] = 0.7778852317840375
$r7[0,2] = 0.7778852317840375
$r7[0,3] = 0.9802868644774461
$r7[0,4] = 0.9802868644774461
$r7[0,5] = 0.7582973271237857
$r7[0,6] = 4.293263378292326
$r7[0,7] = 0.9685679162268007
$r7[0,8] = 0.9802868644774461
$r7[0,9] = 1.044529099493631
$r7[0,10] = 0.9112071656388342
$r7[0,11] = 0.9112071656388342
$r7[0,12] = 0.8602372194671514
$r7[0,13] = 0.9342337319183714
$r7[0,14] = 0.9342337319183714
$r7[0,15] = 0.9457470150581401
$r7[0,16] = 0.9457470150581401
$r7[0,17] = 1.470471636233005
$ws.Range("C7:T7").Value = $r7

# Row 8: HKL index 6 = "Tilt Rotate"
$ws.Cells.Item(8,1).Value = 6
$ws.Cells.Item(8,2).Value = "Tilt Rotate"
$r8 = New-Object 'object[,]' 1,18
$r8[0,0] = 0.5101021293878658
$r8[0,1] = 0.3941785473698862
$r8[0,2] = 0.3941785473698862
$r8[0,3] = 0.4373357528983553
$r8[0,4] = 0.4373357528983553
$r8[0,5] = 0.4052226872482287
$r8[0,6] = 24.38632783153905
$r8[0,7] = 0.4773246116793712
$r8[0,8] = 0.4373357528983553
$r8[0,9] = 0.5101021293878658
$r8[0,10] = 0.452140338378876
$r8[0,11] = 0.452140338378876
$r8[0,12] = 0.4365011213353269
$r8[0,13] = 0.4472054765520357
$r8[0,14] = 0.4472054765520357
$r8[0,15] = 0.4447380456386156
$r8[0,16] = 0.4447380456386156
$r8[0,17] = 4.435081926687126
$ws.Range("C8:T8").Value = $r8

# Row 9: HKL index 7 = "CLR"
$ws.Cells.Item(9,1).Value = 7
$ws.Cells.Item(9,2).Value = "CLR"
$r9 = New-Object 'object[,]' 1,18
$r9[0,0] = 1.081770562969901
$r9[0,1] = 0.9535662553865667
$r9[0,2] = 0.9535662553865667
$r9[0,3] = 1.195262332867592
$r9[0,4] = 1.195262332867592
$r9[0,5] = 0.9521094646260612
$r9[0,6] = 0.9752708429192354
$r9[0,7] = 1.068697454826213
$r9[0,8] = 1.195262332867592
$r9[0,9] = 1.081770562969901
$r9[0,10] = 1.017668409178234
$r9[0,11] = 1.017668409178234
$r9[0,12] = 0.9958154276608431
$r9[0,13] = 1.076866383741353
$r9[0,14] = 1.076866383741353
$r9[0,15] = 1.106465371022913
$r9[0,16] = 1.106465371022913
$r9[0,17] = 1.037779485599262
$ws.Range("C9:T9").Value = $r9

# Row 10: HKL index 8 = "Rizzie Hex"
$ws.Cells.Item(10,1).Value = 8
$ws.Cells.Item(10,2).Value = "Rizzie Hex"
$r10 = New-Object 'object[,]' 1,18
$r10[0,0] = 1.028968515598282
$r10[0,1] = 1.023588681706443
$r10[0,2] = 1.023588681706443
$r10[0,3] = 0.8986392809310543
$r10[0,4] = 0.8986392809310543
$r10[0,5] = 0.9954442626852554
$r10[0,6] = 1.043588727727399
$r10[0,7] = 1.006814243809744
$r10[0,8] = 0.8986392809310543
$r10[0,9] = 1.028968515598282
$r10[0,10] = 1.026278598652362
$r10[0,11] = 1.026278598652362
$r10[0,12] = 1.016000486663327
$r10[0,13] = 0.9837321594119263
$r10[0,14] = 0.9837321594119265
$r10[0,15] = 0.9624589397917084
$r10[0,16] = 0.9624589397917084
$r10[0,17] = 0.9995072854096961
$ws.Range("C10:T10").Value = $r10

# Row 11: HKL index 9 = "Matthies Hex"
$ws.Cells.Item(11,1).Value = 9
$ws.Cells.Item(11,2).Value = "Matthies Hex"
$r11 = New-Object 'object[,]' 1,18
$r11[0,0] = 1.187189153562754
$r11[0,1] = 0.9926561183632608
$r11[0,2] = 0.9926561183632608
$r11[0,3] = 1.425935003305045
$r11[0,4] = 1.425935003305045
$r11[0,5] = 0.8893813698468173
$r11[0,6] = 1.036918289223419
$r11[0,7] = 1.166803114073586
$r11[0,8] = 1.425935003305045
$r11[0,9] = 1.187189153562754
$r11[0,10] = 1.089922635963008
$r11[0,11] = 1.089922635963008
$r11[0,12] = 1.023075547257611
$r11[0,13] = 1.201926758410353
$r11[0,14] = 1.201926758410353
$r11[0,15] = 1.257928819634026
$r11[0,16] = 1.257928819634026
$r11[0,17] = 1.11648050806248
$ws.Range("C11:T11").Value = $r11

# Row 12: HKL index 10 = "Tilt Rotate_Partial"
$ws.Cells.Item(12,1).Value = 10
$ws.Cells.Item(12,2).Value = "Tilt Rotate_Partial"
$r12 = New-Object 'object[,]' 1,18
$r12[0,0] = 0.4687951323330179
$r12[0,1] = 0.5131736293594477
$r12[0,2] = 0.5131736293594477
$r12[0,3] = 0.4158186291731588
$r12[0,4] = 0.4158186291731588
$r12[0,5] = 0.3593355396541261
$r12[0,6] = 24.80948701685619
$r12[0,7] = 0.4665504570366799
$r12[0,8] = 0.4158186291731588
$r12[0,9] = 0.4687951323330179
$r12[0,10] = 0.4909843808462328
$r12[0,11] = 0.4909843808462328
$r12[0,12] = 0.4471014337821972
$r12[0,13] = 0.4659291302885415
$r12[0,14] = 0.4659291302885415
$r12[0,15] = 0.4534015050096958
$r12[0,16] = 0.4534015050096958
$r12[0,17] = 4.505526734068769
$ws.Range("C12:T12").Value = $r12

# Row 13: HKL index 11 = "RotRing OmegaMax-60"
$ws.Cells.Item(13,1).Value = 11
$ws.Cells.Item(13,2).Value = "RotRing OmegaMax-60"
$r13 = New-Object 'object[,]' 1,18
$r13[0,0] = 0.8156180639210705
$r13[0,1] = 1.060154353992671
$r13[0,2] = 1.060154353992671
$r13[0,3] = 1.842711796450164
$r13[0,4] = 1.842711796450164
$r13[0,5] = 0.9351783378329415
$r13[0,6] = 1.040294988388148
$r13[0,7] = 0.899798945790095
$r13[0,8] = 1.842711796450164
$r13[0,9] = 0.8156180639210705
$r13[0,10] = 0.9378862089568709
$r13[0,11] = 0.9378862089568709
$r13[0,12] = 0.9369835852488944
$r13[0,13] = 1.239494738121302
$r13[0,14] = 1.239494738121302
$r13[0,15] = 1.390299002703518
$r13[0,16] = 1.390299002703518
$r13[0,17] = 1.098959414395848
$ws.Range("C13:T13").Value = $r13

# Row 14: HKL index 12 = "Equal Angle_Partial"
$ws.Cells.Item(14,1).Value = 12
$ws.Cells.Item(14,2).Value = "Equal Angle_Partial"
$r14 = New-Object 'object[,]' 1,18
$r14[0,0] = 0.9690302955545643
$r14[0,1] = 1.10858716510821
$r14[0,2] = 1.10858716510821
$r14[0,3] = 0.7924229677633223
$r14[0,4] = 0.7924229677633223
$r14[0,5] = 0.670501355201496
$r14[0,6] = 5.545783250667502
$r14[0,7] = 0.8663474969016852
$r14[0,8] = 0.7924229677633223
$r14[0,9] = 0.9690302955545643
$r14[0,10] = 1.038808730331387
$r14[0,11] = 1.038808730331387
$r14[0,12] = 0.9160396052880903
$r14[0,13] = 0.956680142808699
$r14[0,14] = 0.956680142808699
$r14[0,15] = 0.9156158490473547
$r14[0,16] = 0.9156158490473547
$r14[0,17] = 1.658778755199463
$ws.Range("C14:T14").Value = $r14

# Row 15: HKL index 13 = "Rizzie Hex_Partial"
$ws.Cells.Item(15,1).Value = 13
$ws.Cells.Item(15,2).Value = "Rizzie Hex_Partial"
$r15 = New-Object 'object[,]' 1,18
$r15[0,0] = 0.8351927237057024
$r15[0,1] = 1.678372903941501
$r15[0,2] = 1.678372903941501
$r15[0,3] = 1.01728161075195
$r15[0,4] = 1.01728161075195
$r15[0,5] = 1.021635567048768
$r15[0,6] = 0.5627965227288613
$r15[0,7] = 0.9726220859094149
$r15[0,8] = 1.01728161075195
$r15[0,9] = 0.8351927237057024
$r15[0,10] = 1.256782813823602
$r15[0,11] = 1.256782813823602
$r15[0,12] = 1.178400398231991
$r15[0,13] = 1.176949079466385
$r15[0,14] = 1.176949079466385
$r15[0,15] = 1.137032212287776
$r15[0,16] = 1.137032212287776
$r15[0,17] = 1.014650235681033
$ws.Range("C15:T15").Value = $r15

# Row 16: HKL index 14 = "ND Single"
$ws.Cells.Item(16,1).Value = 14
$ws.Cells.Item(16,2).Value = "ND Single"
$r16 = New-Object 'object[,]' 1,18
$r16[0,0] = 0.0080378142
$r16[0,1] = 0.006236067899999991
$r16[0,2] = 0.006236067899999991
$r16[0,3] = -0.002792170799999996
$r16[0,4] = -0.002792170799999996
$r16[0,5] = -0.002063599
$r16[0,6] = 45.56474799999992
$r16[0,7] = -0.005504551799999996
$r16[0,8] = -0.002792170799999996
$r16[0,9] = 0.0080378142
$r16[0,10] = 0.007136941049999995
$r16[0,11] = 0.007136941049999995
$r16[0,12] = 0.004070094366666664
$r16[0,13] = 0.003827237099999999
$r16[0,14] = 0.003827237099999999
$r16[0,15] = 0.002172385125
$r16[0,16] = 0.002172385125
$r16[0,17] = 7.594776926749987
$ws.Range("C16:T16").Value = $r16

# Row 17: HKL index 15 = "RD Single"
$ws.Cells.Item(17,1).Value = 15
$ws.Cells.Item(17,2).Value = "RD Single"
$r17 = New-Object 'object[,]' 1,18
$r17[0,0] = 15.189881
$r17[0,1] = -0.0064341065
$r17[0,2] = -0.0064341065
$r17[0,3] = -0.0012290963
$r17[0,4] = -0.0012290963
$r17[0,5] = 0.002121282
$r17[0,6] = 0.008039045599999999
$r17[0,7] = 0.56694817
$r17[0,8] = -0.0012290963
$r17[0,9] = 15.189881
$r17[0,10] = 7.59172344675
$r17[0,11] = 7.59172344675
$r17[0,12] = 5.061856058499999
$r17[0,13] = 5.060739265733333
$r17[0,14] = 5.060739265733333
$r17[0,15] = 3.795247175225
$r17[0,16] = 3.795247175225
$r17[0,17] = 2.626554382466666
$ws.Range("C17:T17").Value = $r17

# Row 18: HKL index 16 = "TD Single"
$ws.Cells.Item(18,1).Value = 16
$ws.Cells.Item(18,2).Value = "TD Single"
$r18 = New-Object 'object[,]' 1,18
$r18[0,0] = -0.0012295708
$r18[0,1] = 0.0030969418
$r18[0,2] = 0.0030969418
$r18[0,3] = 30.381976
$r18[0,4] = 30.381976
$r18[0,5] = 0.0025420956
$r18[0,6] = -0.0027921709
$r18[0,7] = 0.0037430372
$r18[0,8] = 30.381976
$r18[0,9] = -0.0012295708
$r18[0,10] = 0.0009336855
$r18[0,11] = 0.0009336855
$r18[0,12] = 0.0014698222
$r18[0,13] = 10.12794779033333
$r18[0,14] = 10.12794779033333
$r18[0,15] = 15.19145484275
$r18[0,16] = 15.19145484275
$r18[0,17] = 5.064556055483332
$ws.Range("C18:T18").Value = $r18

# Row 19: HKL index 17 = "Morris Single"
$ws.Cells.Item(19,1).Value = 17
$ws.Cells.Item(19,2).Value = "Morris Single"
$r19 = New-Object 'object[,]' 1,18
$r19[0,0] = 1.1185846
$r19[0,1] = 0.003165176899999999
$r19[0,2] = 0.003165176899999999
$r19[0,3] = -0.0014782426
$r19[0,4] = -0.0014782426
$r19[0,5] = 0.011519388
$r19[0,6] = -0.0057362253
$r19[0,7] = 0.05500227199999999
$r19[0,8] = -0.0014782426
$r19[0,9] = 1.1185846
$r19[0,10] = 0.56087488845
$r19[0,11] = 0.56087488845
$r19[0,12] = 0.3777563883
$r19[0,13] = 0.3734238447666667
$r19[0,14] = 0.3734238447666667
$r19[0,15] = 0.279698322925
$r19[0,16] = 0.279698322925
$r19[0,17] = 0.1968428281666667
$ws.Range("C19:T19").Value = $r19

# Row 20: HKL index 18 = "Ring Perpendicular to ND"
$ws.Cells.Item(20,1).Value = 18
$ws.Cells.Item(20,2).Value = "Ring Perpendicular to ND"
$r20 = New-Object 'object[,]' 1,18
$r20[0,0] = 2.866815715448219
$r20[0,1] = -0.001348055636986302
$r20[0,2] = -0.001348055636986302
$r20[0,3] = 5.319488004905479
$r20[0,4] = 5.319488004905479
$r20[0,5] = -0.0003479084849315068
$r20[0,6] = 0.001346627909589041
$r20[0,7] = 2.66879592160548
$r20[0,8] = 5.319488004905479
$r20[0,9] = 2.866815715448219
$r20[0,10] = 1.432733829905616
$r20[0,11] = 1.432733829905616
$r20[0,12] = 0.9550399171087672
$r20[0,13] = 2.728318554905571
$r20[0,14] = 2.728318554905571
$r20[0,15] = 3.376110917405548
$r20[0,16] = 3.376110917405548
$r20[0,17] = 1.809125050957808
$ws.Range("C20:T20").Value = $r20

# Row 21: HKL index 19 = "Ring Perpendicular to RD"
$ws.Cells.Item(21,1).Value = 19
$ws.Cells.Item(21,2).Value = "Ring Perpendicular to RD"
$r21 = New-Object 'object[,]' 1,18
$r21[0,0] = 0.2667974827382105
$r21[0,1] = 0.001263986031999999
$r21[0,2] = 0.001263986031999999
$r21[0,3] = 2.504722204034894
$r21[0,4] = 2.504722204034894
$r21[0,5] = 0.3741772924495264
$r21[0,6] = 3.752783920057364
$r21[0,7] = 1.148879533968421
$r21[0,8] = 2.504722204034894
$r21[0,9] = 0.2667974827382105
$r21[0,10] = 0.1340307343851053
$r21[0,11] = 0.1340307343851053
$r21[0,12] = 0.2140795870732456
$r21[0,13] = 0.9242612242683683
$r21[0,14] = 0.9242612242683683
$r21[0,15] = 1.31937646921
$r21[0,16] = 1.31937646921
$r21[0,17] = 1.341437403213403
$ws.Range("C21:T21").Value = $r21

# Row 22: HKL index 20 = "Ring Perpendicular to TD"
$ws.Cells.Item(22,1).Value = 20
$ws.Cells.Item(22,2).Value = "Ring Perpendicular to TD"
$r22 = New-Object 'object[,]' 1,18
$r22[0,0] = 2.105215038871263
$r22[0,1] = 3.406934148036842
$r22[0,2] = 3.406934148036842
$r22[0,3] = 1.704390385821579
$r22[0,4] = 1.704390385821579
$r22[0,5] = 0.03203716434815789
$r22[0,6] = 6.308808287847365
$r22[0,7] = 0.1417101757157894
$r22[0,8] = 1.704390385821579
$r22[0,9] = 2.105215038871263
$r22[0,10] = 2.756074593454052
$r22[0,11] = 2.756074593454052
$r22[0,12] = 1.848062117085421
$r22[0,13] = 2.405513190909895
$r22[0,14] = 2.405513190909895
$r22[0,15] = 2.230232489637816
$r22[0,16] = 2.230232489637816
$r22[0,17] = 2.283182533440166
$ws.Range("C22:T22").Value = $r22

# Row 23: HKL index 21 = "OffsetFTD"
$ws.Cells.Item(23,1).Value = 21
$ws.Cells.Item(23,2).Value = "OffsetFTD"
$r23 = New-Object 'object[,]' 1,18
$r23[0,0] = 0.3279806751991264
$r23[0,1] = 0.009564223208535151
$r23[0,2] = 0.009564223208535151
$r23[0,3] = 1.824526065424647
$r23[0,4] = 1.824526065424647
$r23[0,5] = 2.033072525764183
$r23[0,6] = 0.0009583326907496712
$r23[0,7] = 1.123644033566309
$r23[0,8] = 1.824526065424647
$r23[0,9] = 0.3279806751991264
$r23[0,10] = 0.1687724492038308
$r23[0,11] = 0.1687724492038308
$r23[0,12] = 0.7902058080572815
$r23[0,13] = 0.7206903212774362
$r23[0,14] = 0.7206903212774362
$r23[0,15] = 0.9966492573142389
$r23[0,16] = 0.9966492573142389
$r23[0,17] = 0.886624309308925
$ws.Range("C23:T23").Value = $r23

# Row 24: HKL index 22 = "OffsetATD"
$ws.Cells.Item(24,1).Value = 22
$ws.Cells.Item(24,2).Value = "OffsetATD"
$r24 = New-Object 'object[,]' 1,18
$r24[0,0] = 1.365474223309173
$r24[0,1] = 3.037575828329734
$r24[0,2] = 3.037575828329734
$r24[0,3] = 0.3692450097986277
$r24[0,4] = 0.3692450097986277
$r24[0,5] = 0.5177535224561548
$r24[0,6] = 0.3699524454777166
$r24[0,7] = 0.9633165888619176
$r24[0,8] = 0.3692450097986277
$r24[0,9] = 1.365474223309173
$r24[0,10] = 2.201525025819454
$r24[0,11] = 2.201525025819454
$r24[0,12] = 1.640267858031687
$r24[0,13] = 1.590765020479179
$r24[0,14] = 1.590765020479179
$r24[0,15] = 1.285385017809041
$r24[0,16] = 1.285385017809041
$r24[0,17] = 1.103886269705554
$ws.Range("C24:T24").Value = $r24

# Row 25: HKL index 23 = "OffsetF45"
$ws.Cells.Item(25,1).Value = 23
$ws.Cells.Item(25,2).Value = "OffsetF45"
$r25 = New-Object 'object[,]' 1,18
$r25[0,0] = 0.5353173396279407
$r25[0,1] = 0.2731877001756435
$r25[0,2] = 0.2731877001756435
$r25[0,3] = 1.839099416743353
$r25[0,4] = 1.839099416743353
$r25[0,5] = 1.28050112250611
$r25[0,6] = 0.001005896292277912
$r25[0,7] = 1.076101711066956
$r25[0,8] = 1.839099416743353
$r25[0,9] = 0.5353173396279407
$r25[0,10] = 0.4042525199017921
$r25[0,11] = 0.4042525199017921
$r25[0,12] = 0.6963353874365646
$r25[0,13] = 0.8825348188489791
$r25[0,14] = 0.8825348188489791
$r25[0,15] = 1.121675968322573
$r25[0,16] = 1.121675968322573
$r25[0,17] = 0.8342021977353801
$ws.Range("C25:T25").Value = $r25

# Row 26: HKL index 24 = "OffsetA45"
$ws.Cells.Item(26,1).Value = 24
$ws.Cells.Item(26,2).Value = "OffsetA45"
$r26 = New-Object 'object[,]' 1,18
$r26[0,0] = 1.255728795392275
$r26[0,1] = 1.598615370252646
$r26[0,2] = 1.598615370252646
$r26[0,3] = 0.9996281589042396
$r26[0,4] = 0.9996281589042396
$r26[0,5] = 0.6970805916650487
$r26[0,6] = 0.4158354810641348
$r26[0,7] = 1.063608608492982
$r26[0,8] = 0.9996281589042396
$r26[0,9] = 1.255728795392275
$r26[0,10] = 1.427172082822461
$r26[0,11] = 1.427172082822461
$r26[0,12] = 1.183808252436657
$r26[0,13] = 1.284657441516387
$r26[0,14] = 1.284657441516387
$r26[0,15] = 1.21340012086335
$r26[0,16] = 1.21340012086335
$r26[0,17] = 1.005082834295221
$ws.Range("C26:T26").Value = $r26

# Row 27: HKL index 25 = "OffsetFRD"
$ws.Cells.Item(27,1).Value = 25
$ws.Cells.Item(27,2).Value = "OffsetFRD"
$r27 = New-Object 'object[,]' 1,18
$r27[0,0] = 0.6840436776613265
$r27[0,1] = 0.9200804871536474
$r27[0,2] = 0.9200804871536474
$r27[0,3] = 0.3628513711394295
$r27[0,4] = 0.3628513711394295
$r27[0,5] = 1.655224724001139
$r27[0,6] = 0.0007741736952879507
$r27[0,7] = 0.7786431169962688
$r27[0,8] = 0.3628513711394295
$r27[0,9] = 0.6840436776613265
$r27[0,10] = 0.8020620824074869
$r27[0,11] = 0.8020620824074869
$r27[0,12] = 1.086449629605371
$r27[0,13] = 0.6556585119848011
$r27[0,14] = 0.6556585119848011
$r27[0,15] = 0.5824567267734582
$r27[0,16] = 0.5824567267734582
$r27[0,17] = 0.7336029251078497
$ws.Range("C27:T27").Value = $r27

# Row 28: HKL index 26 = "OffsetARD"
$ws.Cells.Item(28,1).Value = 26
$ws.Cells.Item(28,2).Value = "OffsetARD"
$r28 = New-Object 'object[,]' 1,18
$r28[0,0] = 1.219085589634244
$r28[0,1] = 0.9183225674123235
$r28[0,2] = 0.9183225674123235
$r28[0,3] = 0.3944422404019887
$r28[0,4] = 0.3944422404019887
$r28[0,5] = 1.105831862761006
$r28[0,6] = 0.4993208193803643
$r28[0,7] = 0.6052202324906033
$r28[0,8] = 0.3944422404019887
$r28[0,9] = 1.219085589634244
$r28[0,10] = 1.068704078523284
$r28[0,11] = 1.068704078523284
$r28[0,12] = 1.081080006602524
$r28[0,13] = 0.8439501324828521
$r28[0,14] = 0.8439501324828521
$r28[0,15] = 0.7315731594626362
$r28[0,16] = 0.7315731594626362
$r28[0,17] = 0.7903705520134215
$ws.Range("C28:T28").Value = $r28

# Row 29: HKL index 27 = "Gaussian Quadrature"
$ws.Cells.Item(29,1).Value = 27
$ws.Cells.Item(29,2).Value = "Gaussian Quadrature"
$r29 = New-Object 'object[,]' 1,18
$r29[0,0] = 0.2736435447542511
$r29[0,1] = 0.006185247811692524
$r29[0,2] = 0.006185247811692524
$r29[0,3] = 1.234476484103664
$r29[0,4] = 1.234476484103664
$r29[0,5] = 1.684837296811535
$r29[0,6] = 5.329757293588385
$r29[0,7] = 0.649220902073655
$r29[0,8] = 1.234476484103664
$r29[0,9] = 0.2736435447542511
$r29[0,10] = 0.1399143962829718
$r29[0,11] = 0.1399143962829718
$r29[0,12] = 0.6548886964591597
$r29[0,13] = 0.5047684255565358
$r29[0,14] = 0.5047684255565359
$r29[0,15] = 0.6871954401933179
$r29[0,16] = 0.6871954401933179
$r29[0,17] = 1.529686794857197
$ws.Range("C29:T29").Value = $r29

# Row 30: HKL index 28 = "Michael-CCHex"
$ws.Cells.Item(30,1).Value = 28
$ws.Cells.Item(30,2).Value = "Michael-CCHex"
$r30 = New-Object 'object[,]' 1,18
$r30[0,0] = 0.1429945048897298
$r30[0,1] = 1.393167664550069
$r30[0,2] = 1.393167664550069
$r30[0,3] = 2.097950431515033
$r30[0,4] = 2.097950431515033
$r30[0,5] = 0.06195448557607755
$r30[0,6] = 1.454077748140057
$r30[0,7] = 0.4047026339123727
$r30[0,8] = 2.097950431515033
$r30[0,9] = 0.1429945048897298
$r30[0,10] = 0.7680810847198994
$r30[0,11] = 0.7680810847198994
$r30[0,12] = 0.5327055516719589
$r30[0,13] = 1.211370866984944
$r30[0,14] = 1.211370866984944
$r30[0,15] = 1.433015758117466
$r30[0,16] = 1.433015758117466
$r30[0,17] = 0.9258079114305565
$ws.Range("C30:T30").Value = $r30

# Row 31: HKL index 29 = "Michael-SNHex"
$ws.Cells.Item(31,1).Value = 29
$ws.Cells.Item(31,2).Value = "Michael-SNHex"
$r31 = New-Object 'object[,]' 1,18
$r31[0,0] = 0.3168520110703535
$r31[0,1] = 0.004659716106946445
$r31[0,2] = 0.004659716106946445
$r31[0,3] = 3.813072930177159
$r31[0,4] = 3.813072930177159
$r31[0,5] = 0.09956031504369223
$r31[0,6] = -0.001530820210306131
$r31[0,7] = 0.04351720067976599
$r31[0,8] = 3.813072930177159
$r31[0,9] = 0.3168520110703535
$r31[0,10] = 0.16075586358865
$r31[0,11] = 0.16075586358865
$r31[0,12] = 0.1403573474069974
$r31[0,13] = 1.37819488578482
$r31[0,14] = 1.37819488578482
$r31[0,15] = 1.986914396882904
$r31[0,16] = 1.986914396882904
$r31[0,17] = 0.7126885588112685
$ws.Range("C31:T31").Value = $r31

